$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 and J1 - copy style from H1 (bold header style) then set values
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:J47
$data = @(
  @(6,7),
  @(6,6),
  @(3,3),
  @(7,7),
  @(8,8),
  @(6,6),
  @(10,10),
  @(7,7),
  @(10,10),
  @(6,7),
  @(9,10),
  @(10,10),
  @(1,2),
  @(6,6),
  @(4,4),
  @(7,7),
  @(1,1),
  @(6,7),
  @(9,9),
  @(7,8),
  @(5,6),
  @(1,2),
  @(1,2),
  @(6,7),
  @(1,2),
  @(1,2),
  @(7,7),
  @(3,3),
  @(8,9),
  @(6,7),
  @(7,9),
  @(8,8),
  @(1,2),
  @(7,8),
  @(9,9),
  @(5,6),
  @(5,6),
  @(5,6),
  @(4,5),
  @(5,6),
  @(8,9),
  @(7,8),
  @(5,6),
  @(6,6),
  @(4,4),
  @(1,1)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = 2 + $i
  $ws.Cells.Item($row, 9).Value = $data[$i][0]
  $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
